# Update the "想去人数" (column F) counts on the 展览, 演出 and 全部类型
# sheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value  = 562
$ws1.Cells.Item(3, 6).Value  = 5269
$ws1.Cells.Item(8, 6).Value  = 363
$ws1.Cells.Item(9, 6).Value  = 1332
$ws1.Cells.Item(12, 6).Value = 3050
$ws1.Cells.Item(13, 6).Value = 1887
$ws1.Cells.Item(16, 6).Value = 183
$ws1.Cells.Item(17, 6).Value = 122
$ws1.Cells.Item(18, 6).Value = 651
$ws1.Cells.Item(22, 6).Value = 3472
$ws1.Cells.Item(23, 6).Value = 1080
$ws1.Cells.Item(24, 6).Value = 2748
$ws1.Cells.Item(25, 6).Value = 278
$ws1.Cells.Item(26, 6).Value = 1724
$ws1.Cells.Item(27, 6).Value = 3959
$ws1.Cells.Item(29, 6).Value = 912
$ws1.Cells.Item(30, 6).Value = 454
$ws1.Cells.Item(31, 6).Value = 1264
$ws1.Cells.Item(33, 6).Value = 979
$ws1.Cells.Item(34, 6).Value = 1242
$ws1.Cells.Item(35, 6).Value = 52
$ws1.Cells.Item(36, 6).Value = 997
$ws1.Cells.Item(37, 6).Value = 648
$ws1.Cells.Item(38, 6).Value = 485
$ws1.Cells.Item(39, 6).Value = 398
$ws1.Cells.Item(41, 6).Value = 3541

# --- Sheet: 演出 (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(11, 6).Value = 31
$ws2.Cells.Item(16, 6).Value = 10
$ws2.Cells.Item(22, 6).Value = 34

# --- Sheet: 全部类型 (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value  = 562
$ws4.Cells.Item(3, 6).Value  = 5269
$ws4.Cells.Item(8, 6).Value  = 363
$ws4.Cells.Item(9, 6).Value  = 1332
$ws4.Cells.Item(10, 6).Value = 3050
$ws4.Cells.Item(12, 6).Value = 1887
$ws4.Cells.Item(16, 6).Value = 183
$ws4.Cells.Item(18, 6).Value = 31
$ws4.Cells.Item(20, 6).Value = 122
$ws4.Cells.Item(24, 6).Value = 3472
$ws4.Cells.Item(26, 6).Value = 10
$ws4.Cells.Item(27, 6).Value = 1080
$ws4.Cells.Item(29, 6).Value = 2748
$ws4.Cells.Item(30, 6).Value = 1724
$ws4.Cells.Item(31, 6).Value = 3959
$ws4.Cells.Item(34, 6).Value = 912
$ws4.Cells.Item(35, 6).Value = 1264
$ws4.Cells.Item(37, 6).Value = 979
$ws4.Cells.Item(39, 6).Value = 1242
$ws4.Cells.Item(40, 6).Value = 52
$ws4.Cells.Item(41, 6).Value = 997
$ws4.Cells.Item(42, 6).Value = 648
$ws4.Cells.Item(43, 6).Value = 398
$ws4.Cells.Item(44, 6).Value = 34
$ws4.Cells.Item(48, 6).Value = 3541
